$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 230, 189.35),
    @(3, 110, 100.25),
    @(4, 130, 149.95),
    @(5, 100, 104.35),
    @(6, 158, 143.18),
    @(7, 155, 150.66),
    @(8, 215, 174.21),
    @(9, 105, 125.99),
    @(10, 169, 156.13),
    @(11, 179, 195.43),
    @(12, 65, 66.95),
    @(13, 120, 114.68),
    @(14, 125, 132.56),
    @(15, 95, 112.65),
    @(16, 115, 131.39),
    @(17, 193, 182.13),
    @(18, 120, 112.76),
    @(19, 165, 184.8),
    @(20, 180, 159.52),
    @(21, 125, 112.96)
)

foreach ($row in $data) {
    $r = $row[0]
    $a = $row[1]
    $b = $row[2]
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
}
